$d = $word.ActiveDocument

# Update the date heading (unique in the document).
$d.Content.Find.Execute("2025-10-01 Wednesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-10-02 Thursday", 2) | Out-Null

# Update each division-problem cell in the first table by explicit
# row/column address so that duplicate problem text (e.g. "34÷7=4, 6"
# or "16÷9=1, 7" appearing more than once) is only changed where the
# diff says it should be.
$t = $d.Tables.Item(1)

$changes = @(
    @(1,1,"44÷4=11, 0"),
    @(1,2,"91÷8=11, 3"),
    @(1,3,"12÷7=1, 5"),
    @(1,4,"93÷3=31, 0"),
    @(1,5,"51÷9=5, 6"),

    @(5,1,"16÷9=1, 7"),
    @(5,2,"25÷9=2, 7"),
    @(5,3,"25÷4=6, 1"),
    @(5,4,"36÷3=12, 0"),
    @(5,5,"75÷7=10, 5"),

    @(9,1,"54÷7=7, 5"),
    @(9,2,"57÷4=14, 1"),
    @(9,3,"21÷4=5, 1"),
    @(9,4,"88÷2=44, 0"),
    @(9,5,"45÷5=9, 0"),

    @(13,1,"33÷5=6, 3"),
    @(13,2,"37÷6=6, 1"),
    @(13,3,"69÷5=13, 4"),
    @(13,4,"82÷4=20, 2"),
    @(13,5,"38÷2=19, 0"),

    @(17,1,"52÷5=10, 2"),
    @(17,2,"64÷6=10, 4"),
    @(17,3,"28÷6=4, 4"),
    @(17,4,"66÷5=13, 1"),
    @(17,5,"96÷5=19, 1")
)

foreach ($change in $changes) {
    $row = $change[0]
    $col = $change[1]
    $newText = $change[2]
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $newText
}
